$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# E2 previously held "Random Sounds" (no longer used anywhere) - it now holds
# what used to be in D5 ("can't stop the feeling"), shifting the "body" column
# sounds up one row.
$ws.Range("E2").Value = "can't stop the feeling"

# D5 gets the new "Open 2 inner doors" entry, with the cell set to wrap text.
$ws.Range("D5").Value = "Open 2 inner doors"
$ws.Range("D5").WrapText = $true

# Update the active selection to E2 (was B5).
$ws.Range("E2").Select()
